$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1462729.5
$ws.Range("J19").Value = 1073.8572
$ws.Range("L19").Value = 1073.8572
$ws.Range("N19").Value = -1423.8572

$ws.Range("H34").Value = 2434.8572
$ws.Range("I34").Value = 2434.8572
$ws.Range("K34").Value = 2434.8572
$ws.Range("M34").Value = -2231.8572

$ws.Range("H36").Value = 2434.8572
$ws.Range("I36").Value = 2434.8572
$ws.Range("K36").Value = 2434.8572
$ws.Range("M36").Value = -1719.8572

$ws.Range("H64").Value = 3900
$ws.Range("I64").Value = 3675
$ws.Range("J64").Value = 4080
$ws.Range("K64").Value = 3675
$ws.Range("L64").Value = 4080
$ws.Range("M64").Value = -3427
$ws.Range("N64").Value = -4576

$ws.Range("H67").Value = 3900
$ws.Range("I67").Value = 3675
$ws.Range("J67").Value = 4080
$ws.Range("K67").Value = 3675
$ws.Range("L67").Value = 4080
$ws.Range("M67").Value = -2817
$ws.Range("N67").Value = -5796

$ws.Range("H93").Value = 24660
$ws.Range("J93").Value = 24660
$ws.Range("L93").Value = 24660
$ws.Range("N93").Value = -29652

$ws.Range("H135").Value = 1160.1111
$ws.Range("I135").Value = 616.5454999999999
$ws.Range("J135").Value = 2014.2858
$ws.Range("K135").Value = 5548.9095
$ws.Range("L135").Value = 18128.5722
$ws.Range("M135").Value = -3013.9095
$ws.Range("N135").Value = -23198.5722

$ws.Range("H141").Value = 335883.66
$ws.Range("I141").Value = 402239.4
$ws.Range("J141").Value = 4105
$ws.Range("K141").Value = 1206718.2
$ws.Range("L141").Value = 12315
$ws.Range("M141").Value = -1201538.2
$ws.Range("N141").Value = -22675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1326
$ws.Range("I45").Value = 1326
$ws.Range("K45").Value = 1326
$ws.Range("M45").Value = -949

$ws.Range("H61").Value = 2148.0952
$ws.Range("J61").Value = 2190
$ws.Range("L61").Value = 2190
$ws.Range("N61").Value = -2614

$ws.Range("H102").Value = 2701.2856
$ws.Range("I102").Value = 2477.5
$ws.Range("J102").Value = 2999.6667
$ws.Range("K102").Value = 2477.5
$ws.Range("L102").Value = 2999.6667
$ws.Range("M102").Value = -855.5
$ws.Range("N102").Value = -6243.6667

$ws.Range("H103").Value = 34586.117
$ws.Range("J103").Value = 34586.117
$ws.Range("L103").Value = 34586.117
$ws.Range("N103").Value = -36930.117

$ws.Range("H136").Value = 2148.0952
$ws.Range("J136").Value = 2190
$ws.Range("L136").Value = 6570
$ws.Range("N136").Value = -11670

$ws.Range("H137").Value = 39148.75
$ws.Range("J137").Value = 40484.285
$ws.Range("L137").Value = 40484.285
$ws.Range("N137").Value = -50684.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 29653.846
$ws.Range("J95").Value = 29653.846
$ws.Range("L95").Value = 29653.846
$ws.Range("N95").Value = -35145.84600000001

$ws.Range("H134").Value = 1736.6333
$ws.Range("I134").Value = 1146.7084
$ws.Range("J134").Value = 4096.3335
$ws.Range("K134").Value = 3440.1252
$ws.Range("L134").Value = 12289.0005
$ws.Range("M134").Value = -905.1251999999999
$ws.Range("N134").Value = -17359.0005

$ws.Range("H137").Value = 45510
$ws.Range("J137").Value = 45510
$ws.Range("L137").Value = 45510
$ws.Range("N137").Value = -55710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3505.2083
$ws.Range("I31").Value = 1468.3334
$ws.Range("K31").Value = 1468.3334
$ws.Range("M31").Value = -1173.3334

$ws.Range("H34").Value = 3505.2083
$ws.Range("I34").Value = 1468.3334
$ws.Range("K34").Value = 1468.3334
$ws.Range("M34").Value = -1266.3334

$ws.Range("H62").Value = 3332.5
$ws.Range("I62").Value = 3599
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 3599
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -2975
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 3332.5
$ws.Range("I65").Value = 3599
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 17995
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -14875
$ws.Range("N65").Value = -16240

$ws.Range("H107").Value = 873.5
$ws.Range("I107").Value = 842
$ws.Range("J107").Value = 905
$ws.Range("K107").Value = 842
$ws.Range("L107").Value = 905
$ws.Range("M107").Value = 1078
$ws.Range("N107").Value = -4745

$ws.Range("H137").Value = 41447.5
$ws.Range("J137").Value = 41447.5
$ws.Range("L137").Value = 41447.5
$ws.Range("N137").Value = -51647.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 2615.3333
$ws.Range("J100").Value = 2615.3333
$ws.Range("L100").Value = 7845.999899999999
$ws.Range("N100").Value = -9467.999899999999

$ws.Range("H113").Value = 3906841.2
$ws.Range("I113").Value = 599.63635
$ws.Range("J113").Value = 12500573
$ws.Range("K113").Value = 1798.90905
$ws.Range("L113").Value = 37501719
$ws.Range("M113").Value = 371.09095
$ws.Range("N113").Value = -37506059

$ws.Range("H129").Value = 2076
$ws.Range("I129").Value = 2573.3333
$ws.Range("J129").Value = 1777.6
$ws.Range("K129").Value = 7719.999899999999
$ws.Range("L129").Value = 5332.799999999999
$ws.Range("M129").Value = -2719.999899999999
$ws.Range("N129").Value = -15332.8

$ws.Range("H131").Value = 673.51514
$ws.Range("J131").Value = 803.013
$ws.Range("L131").Value = 2409.039
$ws.Range("N131").Value = -12489.039

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1955
$ws.Range("I97").Value = 1955
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1955
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1459
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 4925.4287
$ws.Range("I102").Value = 2375
$ws.Range("J102").Value = 8326
$ws.Range("K102").Value = 2375
$ws.Range("L102").Value = 8326
$ws.Range("M102").Value = -753
$ws.Range("N102").Value = -11570

$ws.Range("H132").Value = 2776.2856
$ws.Range("I132").Value = 2043.15
$ws.Range("J132").Value = 6034.6665
$ws.Range("K132").Value = 6129.450000000001
$ws.Range("L132").Value = 18103.9995
$ws.Range("M132").Value = -3599.450000000001
$ws.Range("N132").Value = -23163.9995

$ws.Range("H137").Value = 40460
$ws.Range("J137").Value = 40460
$ws.Range("L137").Value = 40460
$ws.Range("N137").Value = -50660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 40400

$ws.Range("H90").Value = 40400

$ws.Range("H132").Value = 5337.3228
$ws.Range("I132").Value = 2425.889
$ws.Range("K132").Value = 7277.667
$ws.Range("M132").Value = -4747.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 72289.89
$ws.Range("J46").Value = 72289.89
$ws.Range("L46").Value = 72289.89
$ws.Range("N46").Value = -72751.89

$ws.Range("H122").Value = 11249.75
$ws.Range("J122").Value = 11249.75
$ws.Range("L122").Value = 33749.25
$ws.Range("N122").Value = -38649.25

$ws.Range("H132").Value = 7579918
$ws.Range("I132").Value = 3989.8235
$ws.Range("J132").Value = 33338074
$ws.Range("K132").Value = 11969.4705
$ws.Range("L132").Value = 100014222
$ws.Range("M132").Value = -9439.470499999999
$ws.Range("N132").Value = -100019282

$ws.Range("H134").Value = 72289.89
$ws.Range("J134").Value = 72289.89
$ws.Range("L134").Value = 216869.67
$ws.Range("N134").Value = -221939.67

$ws.Range("H136").Value = 15678.75
$ws.Range("I136").Value = 13770.875
$ws.Range("J136").Value = 19494.5
$ws.Range("K136").Value = 41312.625
$ws.Range("L136").Value = 58483.5
$ws.Range("M136").Value = -38762.625
$ws.Range("N136").Value = -63583.5
